$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-10 : Name, Email, Phone, Address, Department, Position
$data = @(
    @("admin inv 2", "admininv2@gmail.com", "12346", "pkr", "Civil", "professor"),
    @("admin inv 3", "admininv3@gmail.com", "12347", "ktm", "Computer", "teacher"),
    @("admin inv 4", "admininv4@gmail.com", "12348", "pkr", "Civil", "professor"),
    @("admin inv 7", "admininv7@gmail.com", "12351", "ktm", "Computer", "teacher"),
    @("admin inv 8", "admininv8@gmail.com", "12352", "pkr", "Civil", "professor"),
    @("admin inv 9", "admininv9@gmail.com", "12353", "ktm", "Computer", "teacher"),
    @("admin inv 1", "admininv1@gmail.com", "12345", "ktm", "Computer", "teacher"),
    @("admin inv 5", "-", "12349", "ktm", "Computer", "teacher"),
    @("admin inv 6", "-", "12350", "pkr", "Civil", "professor")
)

# Copy the formatting (bold + border style) already used in column A
# of the existing data rows down onto the new rows that will be added.
$ws.Range("A4").Copy()
$ws.Range("A5:A10").PasteSpecial(-4122)  # xlPasteFormats

# Force the phone number column into text so numeric-looking values
# (e.g. "12346") are preserved exactly as typed, not converted to
# numbers; ClearFormats afterwards drops the temporary number format
# again so the cells end up with the default (no) style, just like
# the rest of the data cells.
$ws.Range("C2:C10").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}

$ws.Range("C2:C10").ClearFormats()
